$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect so we can edit locked (word) cells too.
$ws.Unprotect()

# --- 1. Simple lower-casing fixes to a handful of existing cells ---
$ws.Range("F65").Value  = "breda"
$ws.Range("G93").Value  = "zorginstelling1"
$ws.Range("D100").Value = "beroep1"
$ws.Range("H219").Value = "werkinstelling1"
$ws.Range("C268").Value = "werkinstelling2"

# --- 2. Insert a brand-new 7-row utterance block right before row 296 ---
# (this shifts the existing blocks for utterances #43-#46 down by 7 rows,
# turning them into utterances #44-#47)
$ws.Rows("296:302").Insert()

# Re-apply the correct formatting to the freshly inserted rows: the first
# row of the block ("Utt") uses the yellow-fill style, the other six rows
# use the "unlocked data cell" style. Copy these styles from the
# (now shifted) block immediately below, which still has them intact.
$ws.Range("A303:X303").Copy()
$ws.Range("A296:X296").PasteSpecial(-4122)

$ws.Range("A304:X304").Copy()
$ws.Range("A297:X297").PasteSpecial(-4122)
$ws.Range("A298:X298").PasteSpecial(-4122)
$ws.Range("A299:X299").PasteSpecial(-4122)
$ws.Range("A300:X300").PasteSpecial(-4122)
$ws.Range("A301:X301").PasteSpecial(-4122)
$ws.Range("A302:X302").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 3. Renumber the shifted utterance blocks (old 43-46 -> new 44-47) ---
for ($r = 303; $r -le 330; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 1
}

# --- 4. Populate the new utterance block (#43) at rows 296-302 ---
$ws.Range("A296").Value = 43
$ws.Range("B296").Value = "Utt"
$ws.Range("C296").Value = "ik"
$ws.Range("D296").Value = "weet"

$ws.Range("A297").Value = 43
$ws.Range("B297").Value = "SampleGrootte"

$ws.Range("A298").Value = 43
$ws.Range("B298").Value = "MLU"

$ws.Range("A299").Value = 43
$ws.Range("B299").Value = "Taalmaat"
$ws.Range("D299").Value = "LEX,PV"

$ws.Range("A300").Value = 43
$ws.Range("B300").Value = "Lemma"

$ws.Range("A301").Value = 43
$ws.Range("B301").Value = "Grammaticale fout"

$ws.Range("A302").Value = 43
$ws.Range("B302").Value = "QA"

# Restore protection to match the original workbook's protection state.
$ws.Protect()
